$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$sh = $s.Shapes.Item(1)

# Resize / reposition the "Rectangle 3" text box (grew taller to fit new copy)
$sh.Top = 179.99999
$sh.Height = 277.25

$tr = $sh.TextFrame.TextRange

# "<Team Name>" -> italic "Team Name In Discussion"
$teamRun = $tr.Find("<Team Name>", 0)
$teamRun.Font.Italic = $true
$teamRun.Text = "Team Name In Discussion"

# "Bill Michael" -> "Mr. Bill Michael"
$billRun = $tr.Find("Bill Michael", 0)
$billRun.Text = "Mr. Bill Michael"
